$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 830, pushing existing rows 830-867 down to 832-869.
$ws.Range("A830:A831").EntireRow.Insert()

# New row 830: Red Globe entry dated 2023-04-29 (serial 45041)
$ws.Cells.Item(830, 1).Value2 = 9
$ws.Cells.Item(830, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(830, 3).Value2 = "Metropolitana"
$ws.Cells.Item(830, 4).Value2 = 45041
$ws.Cells.Item(830, 5).Value2 = 13
$ws.Cells.Item(830, 6).Value2 = "Fruta"
$ws.Cells.Item(830, 7).Value2 = 100109
$ws.Cells.Item(830, 8).Value2 = "Uva"
$ws.Cells.Item(830, 9).Value2 = 100109001
$ws.Cells.Item(830, 10).Value2 = "Uva"
$ws.Cells.Item(830, 11).Value2 = "Red Globe"
$ws.Cells.Item(830, 12).Value2 = "Primera"
$ws.Cells.Item(830, 13).Value2 = 220
$ws.Cells.Item(830, 14).Value2 = 11000
$ws.Cells.Item(830, 15).Value2 = 11000
$ws.Cells.Item(830, 16).Value2 = 11000
$ws.Cells.Item(830, 17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(830, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(830, 19).Value2 = 611
$ws.Cells.Item(830, 20).Value2 = 18

# New row 831: Timco entry dated 2023-04-29 (serial 45041)
$ws.Cells.Item(831, 1).Value2 = 9
$ws.Cells.Item(831, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(831, 3).Value2 = "Metropolitana"
$ws.Cells.Item(831, 4).Value2 = 45041
$ws.Cells.Item(831, 5).Value2 = 13
$ws.Cells.Item(831, 6).Value2 = "Fruta"
$ws.Cells.Item(831, 7).Value2 = 100109
$ws.Cells.Item(831, 8).Value2 = "Uva"
$ws.Cells.Item(831, 9).Value2 = 100109001
$ws.Cells.Item(831, 10).Value2 = "Uva"
$ws.Cells.Item(831, 11).Value2 = "Timco"
$ws.Cells.Item(831, 12).Value2 = "Primera"
$ws.Cells.Item(831, 13).Value2 = 300
$ws.Cells.Item(831, 14).Value2 = 12000
$ws.Cells.Item(831, 15).Value2 = 12000
$ws.Cells.Item(831, 16).Value2 = 12000
$ws.Cells.Item(831, 17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(831, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(831, 19).Value2 = 667
$ws.Cells.Item(831, 20).Value2 = 18

Write-Host "Inserted rows 830-831; dimension now extends to row 869"
